$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -8
$ws.Range("F6").Value = -5
$ws.Range("F8").Value = -10
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = -9
$ws.Range("F11").Value = -6
